# lots of work on port-level landings
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 36 got re-labeled: col A now holds the "SAN DIEGO AREA TOTALS" label
# (previously in column C), and column C now reads "Totals" (matching the
# pattern used in the port-subtotal rows above it).
$ws.Range("A36").Value = "SAN DIEGO AREA TOTALS"
$ws.Range("C36").Value = "Totals"

# Column A now needs to be as wide as column C, since it holds the same
# long label text.
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# The active selection moved to A32.
$ws.Range("A32").Select()
